{"js": "const pairs = [\n  [\"2025-11-23 Sunday\", \"2025-11-24 Monday\"],\n  [\"168\u00f79=18, 6\", \"592\u00f75=118, 2\"],\n  [\"313\u00f79=34, 7\", \"816\u00f72=408, 0\"],\n  [\"343\u00f77=49, 0\", \"231\u00f78=28, 7\"],\n  [\"385\u00f73=128, 1\", \"132\u00f73=44, 0\"],\n  [\"262\u00f78=32, 6\", \"145\u00f78=18, 1\"],\n  [\"959\u00f76=159, 5\", \"963\u00f75=192, 3\"],\n  [\"804\u00f72=402, 0\", \"878\u00f79=97, 5\"],\n  [\"756\u00f77=108, 0\", \"738\u00f73=246, 0\"],\n  [\"790\u00f74=197, 2\", \"435\u00f76=72, 3\"],\n  [\"355\u00f78=44, 3\", \"484\u00f79=53, 7\"],\n  [\"431\u00f73=143, 2\", \"786\u00f79=87, 3\"],\n  [\"364\u00f72=182, 0\", \"473\u00f79=52, 5\"],\n  [\"425\u00f76=70, 5\", \"408\u00f79=45, 3\"],\n  [\"827\u00f78=103, 3\", \"599\u00f75=119, 4\"],\n  [\"215\u00f72=107, 1\", \"558\u00f75=111, 3\"],\n  [\"792\u00f78=99, 0\", \"555\u00f74=138, 3\"],\n  [\"187\u00f73=62, 1\", \"227\u00f74=56, 3\"],\n  [\"241\u00f73=80, 1\", \"817\u00f74=204, 1\"],\n  [\"645\u00f76=107, 3\", \"164\u00f74=41, 0\"],\n  [\"929\u00f77=132, 5\", \"608\u00f78=76, 0\"],\n  [\"833\u00f78=104, 1\", \"567\u00f75=113, 2\"],\n  [\"288\u00f79=32, 0\", \"459\u00f75=91, 4\"],\n  [\"797\u00f75=159, 2\", \"676\u00f72=338, 0\"],\n  [\"158\u00f78=19, 6\", \"730\u00f74=182, 2\"],\n  [\"588\u00f73=196, 0\", \"823\u00f72=411, 1\"]\n];\n\nfor (const [before, after] of pairs) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + before);\n  }\n\n  for (const item of results.items) {\n    item.insertText(after, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-11-23 Sunday\", \"2025-11-24 Monday\"),\n    @(\"168\u00f79=18, 6\", \"592\u00f75=118, 2\"),\n    @(\"313\u00f79=34, 7\", \"816\u00f72=408, 0\"),\n    @(\"343\u00f77=49, 0\", \"231\u00f78=28, 7\"),\n    @(\"385\u00f73=128, 1\", \"132\u00f73=44, 0\"),\n    @(\"262\u00f78=32, 6\", \"145\u00f78=18, 1\"),\n    @(\"959\u00f76=159, 5\", \"963\u00f75=192, 3\"),\n    @(\"804\u00f72=402, 0\", \"878\u00f79=97, 5\"),\n    @(\"756\u00f77=108, 0\", \"738\u00f73=246, 0\"),\n    @(\"790\u00f74=197, 2\", \"435\u00f76=72, 3\"),\n    @(\"355\u00f78=44, 3\", \"484\u00f79=53, 7\"),\n    @(\"431\u00f73=143, 2\", \"786\u00f79=87, 3\"),\n    @(\"364\u00f72=182, 0\", \"473\u00f79=52, 5\"),\n    @(\"425\u00f76=70, 5\", \"408\u00f79=45, 3\"),\n    @(\"827\u00f78=103, 3\", \"599\u00f75=119, 4\"),\n    @(\"215\u00f72=107, 1\", \"558\u00f75=111, 3\"),\n    @(\"792\u00f78=99, 0\", \"555\u00f74=138, 3\"),\n    @(\"187\u00f73=62, 1\", \"227\u00f74=56, 3\"),\n    @(\"241\u00f73=80, 1\", \"817\u00f74=204, 1\"),\n    @(\"645\u00f76=107, 3\", \"164\u00f74=41, 0\"),\n    @(\"929\u00f77=132, 5\", \"608\u00f78=76, 0\"),\n    @(\"833\u00f78=104, 1\", \"567\u00f75=113, 2\"),\n    @(\"288\u00f79=32, 0\", \"459\u00f75=91, 4\"),\n    @(\"797\u00f75=159, 2\", \"676\u00f72=338, 0\"),\n    @(\"158\u00f78=19, 6\", \"730\u00f74=182, 2\"),\n    @(\"588\u00f73=196, 0\", \"823\u00f72=411, 1\")\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Execute(\n        $pair[0],   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $pair[1],   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n}\n"}
